# Update TPM-derived NATMI values on the active sheet to reflect the
# newly computed ligand-receptor statistics (rows 2-6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 3.035631
    "H2" = 9.106892999999999
    "I2" = 0.3844053036145329
    "J2" = 0.384405303614533
    "M2" = 0.019
    "N2" = 0.057
    "Q2" = 0.057676989
    "R2" = 0.5190929009999999
    "S2" = 0.3844053036145329
    "T2" = 0.384405303614533

    "I3" = 0.4031610162601934
    "J3" = 0.4031610162601935
    "M3" = 0.019
    "N3" = 0.057
    "Q3" = 0.06049113599999999
    "R3" = 0.544420224
    "S3" = 0.4031610162601934
    "T3" = 0.4031610162601935

    "E4" = 1
    "F4" = 0.3333333333333333
    "G4" = 0.1452593333333333
    "H4" = 0.435778
    "I4" = 0.01839434968638963
    "J4" = 0.01839434968638963
    "M4" = 0.019
    "N4" = 0.057
    "Q4" = 0.002759927333333333
    "R4" = 0.024839346
    "S4" = 0.01839434968638963
    "T4" = 0.01839434968638963

    "G5" = 0.7016706666666667
    "H5" = 2.105012
    "I5" = 0.08885333087500151
    "J5" = 0.08885333087500152
    "M5" = 0.019
    "N5" = 0.057
    "Q5" = 0.01333174266666667
    "R5" = 0.119985684
    "S5" = 0.08885333087500151
    "T5" = 0.08885333087500152

    "E6" = 3
    "F6" = 1
    "G6" = 0.8306490000000001
    "H6" = 2.491947
    "I6" = 0.1051859995638825
    "J6" = 0.1051859995638825
    "M6" = 0.019
    "N6" = 0.057
    "Q6" = 0.057676989
    "R6" = 0.5190929009999999
    "S6" = 0.3844053036145329
    "T6" = 0.384405303614533
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
